# Update sessions
#
# "Generalized pairwise comparisons" (Topic-contributed Sessions, row 9)
# gains a co-organizer: "Arne Bathke" -> "Arne Bathke, Johan Verbeeck"

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Topic-contributed Sessions")

$ws2.Range("B9").Value = "Arne Bathke, Johan Verbeeck"

# The workbook was left with the "Topic-contributed Sessions" tab active
# and cell B9 selected there (previously "Featured Sessions" was active).
$ws2.Activate()
$ws2.Range("B9").Select()
